$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing last row (142) with revised figures for 01-01-2021
$ws.Range("B142").Value = 39368
$ws.Range("C142").Value = 7787
$ws.Range("D142").Value = 31656
$ws.Range("E142").Value = 10937
$ws.Range("F142").Value = 12265
$ws.Range("G142").Value = 38076

# Append new row 143 for 01-04-2021.
# Writing the date-like text directly via .Value triggers Excel's
# auto-date recognition (turns it into a date serial + date-formatted
# style). Route it through a formula-then-paste-as-values round trip so
# it lands as a plain shared string, matching the rest of column A.
$ws.Range("A143").Formula = '="01-04-2021"'
$ws.Range("A143").Copy()
$ws.Range("A143").PasteSpecial(-4163)

$ws.Range("B143").Value = 41667
$ws.Range("C143").Value = 8218
$ws.Range("D143").Value = 33530
$ws.Range("E143").Value = 10700
$ws.Range("F143").Value = 12671
$ws.Range("G143").Value = 39677
